$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 15: add a new customer record
$ws.Range("A15").Value = "宇崝"
$ws.Range("C15").Value = "1FD4-2332-1A6F-054A"
$ws.Range("I15").Value = "忘記何時給的序號，以及是否付費。"

# Update the selected cell to match the author's final cursor position
$ws.Range("F20").Select()
